# Rescale the res_load p_mw results: columns B & D are multiplied by 70
# (fraction -> scaled load, e.g. p.u. * 70 kW) and columns C & E are
# recomputed from the original "C" fraction times 100 and 70 respectively.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 25; $row++) {
    $bVal = $ws.Cells.Item($row, 2).Value()
    $cVal = $ws.Cells.Item($row, 3).Value()

    $ws.Cells.Item($row, 2).Value = $bVal * 70
    $ws.Cells.Item($row, 3).Value = $cVal * 100
    $ws.Cells.Item($row, 4).Value = $bVal * 70
    $ws.Cells.Item($row, 5).Value = $cVal * 70
}
